$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 5.553084769722144

$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 9.295990156953671

$ws.Range("B4").Value = 3.230985683306322
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 6.201049113329182

$ws.Range("B5").Value = 3.230985683306322
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 0.8054896365839992
$ws.Range("E5").Value = 0.496779210170732
$ws.Range("G5").Value = 6.201049113329182

$ws.Range("B6").Value = 0.04763786555579896
$ws.Range("C6").Value = 0.3127903958511391
$ws.Range("D6").Value = 0.1575252929769615
$ws.Range("E6").Value = 0.496779210170732
$ws.Range("G6").Value = 1.014732764554632

$ws.Range("B7").Value = 0.127881588408715
$ws.Range("C7").Value = 1.667794583268128
$ws.Range("D7").Value = 0.1575252929769615
$ws.Range("E7").Value = 0.496779210170732
$ws.Range("G7").Value = 2.449980674824537

$ws.Range("B8").Value = 3.230985683306322
$ws.Range("C8").Value = 1.667794583268128
$ws.Range("D8").Value = 0.8054896365839992
$ws.Range("E8").Value = 0.496779210170732
$ws.Range("G8").Value = 6.201049113329182

$ws.Range("B9").Value = 3.230985683306322
$ws.Range("C9").Value = 1.667794583268128
$ws.Range("D9").Value = 0.8054896365839992
$ws.Range("E9").Value = 0.496779210170732
$ws.Range("G9").Value = 6.201049113329182

$ws.Range("B10").Value = 3.230985683306322
$ws.Range("C10").Value = 1.667794583268128
$ws.Range("D10").Value = 3.900430680208489
$ws.Range("E10").Value = 0.496779210170732
$ws.Range("G10").Value = 9.295990156953671

$ws.Range("B11").Value = 3.230985683306322
$ws.Range("C11").Value = 1.667794583268128
$ws.Range("D11").Value = 0.8054896365839992
$ws.Range("E11").Value = 0.496779210170732
$ws.Range("G11").Value = 6.201049113329182

$ws.Range("B12").Value = 0.127881588408715
$ws.Range("C12").Value = 0.3127903958511391
$ws.Range("D12").Value = 0.8054896365839992
$ws.Range("E12").Value = 0.496779210170732
$ws.Range("G12").Value = 1.742940831014585

$ws.Range("B13").Value = 3.230985683306322
$ws.Range("C13").Value = 1.667794583268128
$ws.Range("D13").Value = 3.900430680208489
$ws.Range("E13").Value = 0.496779210170732
$ws.Range("G13").Value = 9.295990156953671

$ws.Range("B14").Value = 3.230985683306322
$ws.Range("C14").Value = 1.667794583268128
$ws.Range("D14").Value = 3.900430680208489
$ws.Range("E14").Value = 0.496779210170732
$ws.Range("G14").Value = 9.295990156953671

$ws.Range("B15").Value = 0.6753301551942219
$ws.Range("C15").Value = 1.667794583268128
$ws.Range("D15").Value = 0.8054896365839992
$ws.Range("E15").Value = 0.496779210170732
$ws.Range("G15").Value = 3.645393585217082

$ws.Range("B16").Value = 3.230985683306322
$ws.Range("C16").Value = 1.667794583268128
$ws.Range("D16").Value = 0.8054896365839992
$ws.Range("E16").Value = 0.496779210170732
$ws.Range("G16").Value = 6.201049113329182

$ws.Range("B17").Value = 0.6753301551942219
$ws.Range("C17").Value = 1.667794583268128
$ws.Range("D17").Value = 0.1575252929769615
$ws.Range("E17").Value = 0.496779210170732
$ws.Range("G17").Value = 2.997429241610044

$ws.Range("B18").Value = 3.230985683306322
$ws.Range("C18").Value = 1.667794583268128
$ws.Range("D18").Value = 0.8054896365839992
$ws.Range("E18").Value = 0.496779210170732
$ws.Range("G18").Value = 6.201049113329182

$ws.Range("B19").Value = 3.230985683306322
$ws.Range("C19").Value = 1.667794583268128
$ws.Range("D19").Value = 0.8054896365839992
$ws.Range("E19").Value = 0.496779210170732
$ws.Range("G19").Value = 6.201049113329182

$ws.Range("B20").Value = 3.230985683306322
$ws.Range("C20").Value = 1.667794583268128
$ws.Range("D20").Value = 26.21740644021617
$ws.Range("E20").Value = 0.496779210170732
$ws.Range("G20").Value = 31.61296591696135

$ws.Range("B21").Value = 1.459612070389937
$ws.Range("C21").Value = 1.667794583268128
$ws.Range("D21").Value = 26.21740644021617
$ws.Range("E21").Value = 8.660232485948974
$ws.Range("G21").Value = 38.00504557982321

$ws.Range("B22").Value = 3.230985683306322
$ws.Range("C22").Value = 1.667794583268128
$ws.Range("D22").Value = 3.900430680208489
$ws.Range("E22").Value = 0.496779210170732
$ws.Range("G22").Value = 9.295990156953671

$ws.Range("B23").Value = 0.6753301551942219
$ws.Range("C23").Value = 1.667794583268128
$ws.Range("D23").Value = 0.8054896365839992
$ws.Range("E23").Value = 0.496779210170732
$ws.Range("G23").Value = 3.645393585217082

$ws.Range("B24").Value = 3.230985683306322
$ws.Range("C24").Value = 1.667794583268128
$ws.Range("D24").Value = 0.8054896365839992
$ws.Range("E24").Value = 0.496779210170732
$ws.Range("G24").Value = 6.201049113329182

$ws.Range("B25").Value = 1.459612070389937
$ws.Range("C25").Value = 1.667794583268128
$ws.Range("D25").Value = 0.1575252929769615
$ws.Range("E25").Value = 0.496779210170732
$ws.Range("G25").Value = 3.781711156805759

$ws.Range("B26").Value = 0.3048080303191223
$ws.Range("C26").Value = 0.04240448674262143
$ws.Range("D26").Value = 0.8054896365839992
$ws.Range("E26").Value = 0.496779210170732
$ws.Range("G26").Value = 1.649481363816475
